# Scheduled-runner data refresh.
# Re-pulls current marketboard prices (currentAveragePrice / NQ / HQ columns
# H:J) and recomputes the derived Leve-turn-in profit columns (K:N) for each
# job sheet. Values come from an external market-data pull, not in-sheet
# formulas, so the refresh here simply overwrites the stale cells (and
# adds/clears cells where a column crosses to/from zero) to match the new
# snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 149.75
$ws.Range("J5").Value = 147
$ws.Range("L5").Value = 147
$ws.Range("N5").Value = -377
$ws.Range("H8").Value = 80.26667
$ws.Range("I8").Value = 93
$ws.Range("K8").Value = 279
$ws.Range("M8").Value = -140
$ws.Range("H17").Value = 1157.6666
$ws.Range("J17").Value = 1157.6666
$ws.Range("L17").Value = 3472.9998
$ws.Range("N17").Value = -3808.9998
$ws.Range("H31").Value = 59.25
$ws.Range("I31").Value = 59.25
$ws.Range("K31").Value = 177.75
$ws.Range("M31").Value = 52.25
$ws.Range("H70").Value = 1552.6316
$ws.Range("I70").Value = 1468.1818
$ws.Range("J70").Value = 1668.75
$ws.Range("K70").Value = 4404.5454
$ws.Range("L70").Value = 5006.25
$ws.Range("M70").Value = -4134.5454
$ws.Range("N70").Value = -5546.25
$ws.Range("H73").Value = 1552.6316
$ws.Range("I73").Value = 1468.1818
$ws.Range("J73").Value = 1668.75
$ws.Range("K73").Value = 4404.5454
$ws.Range("L73").Value = 5006.25
$ws.Range("M73").Value = -3468.5454
$ws.Range("N73").Value = -6878.25
$ws.Range("H80").Value = 1275
$ws.Range("J80").Value = 1275
$ws.Range("L80").Value = 3825
$ws.Range("N80").Value = -5821
$ws.Range("H83").Value = 1275
$ws.Range("J83").Value = 1275
$ws.Range("L83").Value = 11475
$ws.Range("N83").Value = -21459
$ws.Range("H111").Value = 2187
$ws.Range("I111").Value = 1500
$ws.Range("K111").Value = 4500
$ws.Range("M111").Value = -1433
$ws.Range("H132").Value = 8198.799999999999
$ws.Range("I132").Value = 8198.799999999999
$ws.Range("K132").Value = 24596.4
$ws.Range("M132").Value = -22066.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("H32").Value = 11269.818
$ws.Range("I32").Value = 8396.9
$ws.Range("K32").Value = 8396.9
$ws.Range("M32").Value = -8109.9
$ws.Range("H61").Value = 10463.23
$ws.Range("I61").Value = 12937.667
$ws.Range("K61").Value = 12937.667
$ws.Range("M61").Value = -12725.667
$ws.Range("H88").Value = 1623.625
$ws.Range("I88").Value = 1873.5
$ws.Range("J88").Value = 1373.75
$ws.Range("K88").Value = 1873.5
$ws.Range("L88").Value = 1373.75
$ws.Range("M88").Value = -1467.5
$ws.Range("N88").Value = -2185.75
$ws.Range("H91").Value = 1623.625
$ws.Range("I91").Value = 1873.5
$ws.Range("J91").Value = 1373.75
$ws.Range("K91").Value = 1873.5
$ws.Range("L91").Value = 1373.75
$ws.Range("M91").Value = -469.5
$ws.Range("N91").Value = -4181.75
$ws.Range("H136").Value = 10463.23
$ws.Range("I136").Value = 12937.667
$ws.Range("K136").Value = 38813.001
$ws.Range("M136").Value = -36263.001
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H22").Value = 1025
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -527
$ws.Range("H94").Value = 2051.75
$ws.Range("I94").Value = 1067.3334
$ws.Range("J94").Value = 5005
$ws.Range("K94").Value = 1067.3334
$ws.Range("L94").Value = 5005
$ws.Range("M94").Value = -616.3334
$ws.Range("N94").Value = -5907
$ws.Range("M4").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 24.416666
$ws.Range("I7").Value = 8.25
$ws.Range("K7").Value = 8.25
$ws.Range("M7").Value = 104.75
$ws.Range("H22").Value = 814.4
$ws.Range("J22").Value = 897.5
$ws.Range("L22").Value = 897.5
$ws.Range("N22").Value = -1597.5
$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 15000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -14826
$ws.Range("H35").Value = 3928
$ws.Range("I35").Value = 3022
$ws.Range("J35").Value = 4230
$ws.Range("K35").Value = 3022
$ws.Range("L35").Value = 4230
$ws.Range("M35").Value = -2728
$ws.Range("N35").Value = -4818
$ws.Range("H86").Value = 2495
$ws.Range("J86").Value = 990
$ws.Range("L86").Value = 990
$ws.Range("N86").Value = -3236
$ws.Range("H89").Value = 2495
$ws.Range("J89").Value = 990
$ws.Range("L89").Value = 4950
$ws.Range("N89").Value = -16182
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 890.75
$ws.Range("J44").Value = 1147.1666
$ws.Range("L44").Value = 3441.4998
$ws.Range("N44").Value = -4237.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 49899.332
$ws.Range("I80").Value = 49899
$ws.Range("J80").Value = 49899.5
$ws.Range("K80").Value = 49899
$ws.Range("L80").Value = 49899.5
$ws.Range("M80").Value = -48901
$ws.Range("N80").Value = -51895.5
$ws.Range("H83").Value = 49899.332
$ws.Range("I83").Value = 49899
$ws.Range("J83").Value = 49899.5
$ws.Range("K83").Value = 249495
$ws.Range("L83").Value = 249497.5
$ws.Range("M83").Value = -244503
$ws.Range("N83").Value = -259481.5
$ws.Range("H98").Value = 52976.832
$ws.Range("J98").Value = 52976.832
$ws.Range("L98").Value = 52976.832
$ws.Range("N98").Value = -58966.832
$ws.Range("H126").Value = 2100
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130
$ws.Range("H132").Value = 2689.875
$ws.Range("I132").Value = 2194.5
$ws.Range("J132").Value = 3515.5
$ws.Range("K132").Value = 6583.5
$ws.Range("L132").Value = 10546.5
$ws.Range("M132").Value = -4053.5
$ws.Range("N132").Value = -15606.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1700
$ws.Range("I7").Value = 1700
$ws.Range("K7").Value = 1700
$ws.Range("M7").Value = -1588
$ws.Range("H9").Value = 677.6667
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 516.5
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 516.5
$ws.Range("M9").Value = -776
$ws.Range("N9").Value = -964.5
$ws.Range("H16").Value = 7175
$ws.Range("I16").Value = 7175
$ws.Range("K16").Value = 7175
$ws.Range("M16").Value = -7005
$ws.Range("H22").Value = 871.5
$ws.Range("I22").Value = 744
$ws.Range("K22").Value = 744
$ws.Range("M22").Value = -449
$ws.Range("H27").Value = 871.5
$ws.Range("I27").Value = 744
$ws.Range("K27").Value = 744
$ws.Range("M27").Value = -637
$ws.Range("H39").Value = 5100
$ws.Range("I39").Value = 1200
$ws.Range("K39").Value = 1200
$ws.Range("M39").Value = -740
$ws.Range("H126").Value = 1700
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630
$ws.Range("H132").Value = 5534.1816
$ws.Range("J132").Value = 6199.2
$ws.Range("L132").Value = 18597.6
$ws.Range("N132").Value = -23657.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1817
$ws.Range("I132").Value = 1704
$ws.Range("J132").Value = 2099.5
$ws.Range("K132").Value = 5112
$ws.Range("L132").Value = 6298.5
$ws.Range("M132").Value = -2582
$ws.Range("N132").Value = -11358.5
